# Apply changes described by the commit:
# "Script 1 fix error with deleting everything in temp folder at start of job;
#  added feature for specifying filename for MZmine3 batch xml template file."

$wb = $excel.ActiveWorkbook

$wsJob = $wb.Worksheets.Item("Job to Run")
$wsAll = $wb.Worksheets.Item("All")

# --- Sheet "Job to Run" ---
# Add new column E header
$wsJob.Range("E1").Value = "MZmine3 batch template"

# Update the selected job row (now Anid_HE_TJGIp11_pos_2018 with 3/3 replicates)
$wsJob.Range("A2").Value = "Anid_HE_TJGIp11_pos_2018"
$wsJob.Range("B2").Value = 3
$wsJob.Range("C2").Value = 3
$wsJob.Range("E2").Value = "MZmine3_batch_params_LCMSMS_HE_for_Commandline_2024_8_test_for_Python_workflow.xml"

# widen column B on this sheet (engine snaps ColumnWidth to whole-point
# increments internally, so 22.45 is the input that lands closest to the
# target stored width of 23.28515625 characters)
$wsJob.Columns.Item(2).ColumnWidth = 22.45

# --- Sheet "All" ---
$wsAll.Range("E1").Value = "MZmine3 batch template"
$templateName = "MZmine3_batch_params_LCMSMS_HE_for_Commandline_2024_8_test_for_Python_workflow.xml"
$wsAll.Range("E2").Value = $templateName
$wsAll.Range("E3").Value = $templateName
$wsAll.Range("E4").Value = $templateName
$wsAll.Range("E5").Value = $templateName

# Update selections to match the final saved cursor positions
$wsJob.Range("A16").Select()
$wsAll.Range("E12").Select()
$wsJob.Activate()
